$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

# Row 9 ("MMST 470"): quantity was blank, now lifted 25 units -> Total recalculates to 11750
$ws1.Range("C9").Value = 25

# Row 14 (50 units lifted previously) -> quantity cleared back to blank -> Total recalculates to 0
$ws1.Range("C14").Value = $null

# Row 32: quantity was blank, now 500 units -> Total recalculates to 13985
$ws1.Range("C32").Value = 500

# Row 43 ("CREDIT Lifting"): quantity updated from 104623 to 315964 -> Total recalculates accordingly
$ws1.Range("C43").Value = 315964

# Move the view/cursor to reflect where the user left off editing (last report 01-03-25)
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("H47").Select()
